$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# dob column (C) values: DD/MM/YYYY -> YYYY-MM-DD, kept as TEXT (not auto-converted to a date)
$dates = @{
    2  = "1987-01-04"
    3  = "1978-07-12"
    4  = "1987-09-09"
    5  = "1978-09-01"
    6  = "1967-11-19"
    7  = "1974-11-12"
    8  = "1975-05-06"
    9  = "1976-07-08"
    10 = "1998-10-15"
    11 = "1989-10-12"
    12 = "1987-10-14"
    13 = "1986-05-01"
    14 = "1986-05-07"
    15 = "1980-01-13"
    16 = "1979-01-05"
    17 = "1999-11-16"
    18 = "1990-07-31"
    19 = "1981-03-29"
    20 = "1984-02-18"
}

# Rows 2-4 keep a quote-prefixed text cell (leading apostrophe); rows 5-20 become
# plain text cells without the quote prefix.
$quotePrefixRows = @(2, 3, 4)

foreach ($row in 2..20) {
    $cell = $ws.Range("C$row")
    $cell.NumberFormat = "@"
    $value = $dates[$row]
    if ($quotePrefixRows -contains $row) {
        $cell.Value = "'" + $value
    } else {
        $cell.Value = $value
    }
}
